$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New header labels for the additional "Dominic" data columns (H:M).
#    Setting these creates the 6 new shared-string entries.
# ---------------------------------------------------------------------------
$ws.Range("H1").Value2 = "Butanal-50p"
$ws.Range("I1").Value2 = "Butyroin-50p"
$ws.Range("J1").Value2 = "Octandiol-50p"
$ws.Range("K1").Value2 = "Butanal-20"
$ws.Range("L1").Value2 = "Butyroin-20"
$ws.Range("M1").Value2 = "Octandiol-20"

# ---------------------------------------------------------------------------
# 2) Re-shuffle cell formatting (borders/bold) so that:
#      - H:J take over the formatting the old B:D columns used to have
#      - K:M take over the formatting the old E:G columns used to have
#      - B:G all adopt the "plain" E:G-style look (the stray bold/border
#        that used to sit one row too high, on row 7, now correctly lands
#        on row 8 for every column)
#    Existing formats are first stashed into scratch cells (column Z) so
#    that every later paste always reads from an untouched template,
#    regardless of paste order.
# ---------------------------------------------------------------------------
$stashOps = @(
  @("G2","Z1"),
  @("E7","Z2"),
  @("F2","Z3"),
  @("F7","Z4"),
  @("G7","Z5"),
  @("B2","Z6"),
  @("B8","Z7"),
  @("C2","Z8"),
  @("C8","Z9"),
  @("D8","Z10"),
)
foreach ($op in $stashOps) {
  $ws.Range($op[0]).Copy() | Out-Null
  $ws.Range($op[1]).PasteSpecial(-4122) | Out-Null
}

$applyOps = @(
  @("Z1","B2:B7"),
  @("Z2","B8"),
  @("Z3","C2:C7"),
  @("Z4","C8"),
  @("Z1","D2:D7"),
  @("Z5","D8"),
  @("Z1","E2:E7"),
  @("Z2","E8"),
  @("Z3","F2:F7"),
  @("Z4","F8"),
  @("Z1","G2:G7"),
  @("Z5","G8"),
  @("Z6","H2:H7"),
  @("Z7","H8"),
  @("Z8","I2:I7"),
  @("Z9","I8"),
  @("Z6","J2:J7"),
  @("Z10","J8"),
  @("Z1","K2:K7"),
  @("Z2","K8"),
  @("Z3","L2:L7"),
  @("Z4","L8"),
  @("Z1","M2:M7"),
  @("Z5","M8"),
)
foreach ($op in $applyOps) {
  $ws.Range($op[0]).Copy() | Out-Null
  $ws.Range($op[1]).PasteSpecial(-4122) | Out-Null
}

# Style index "3" (same border look as row-8 style on C/F/L) also needs to
# land on the stray F9 cell right under the table.
$ws.Range("Z4").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4122) | Out-Null

# Drop the scratch column again - it must not show up in the saved sheet.
$ws.Range("Z1:Z10").Clear() | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Write the refreshed measurement values (new fit results for B:G, plus
#    the brand-new Dominic columns H:M) in one shot.
# ---------------------------------------------------------------------------
$values = New-Object 'object[,]' 7,12
$values[0,0] = [double]"124.1761658031088"
$values[0,1] = [double]"15.794709948902915"
$values[0,2] = [double]"4.7571502860114405"
$values[0,3] = [double]"28.994818652849744"
$values[0,4] = [double]"36.017132551848512"
$values[0,5] = [double]"11.694227769110762"
$values[0,6] = [double]"71.550949913644203"
$values[0,7] = [double]"13.744815148782687"
$values[0,8] = [double]"5.4123764950598021"
$values[0,9] = [double]"59.927461139896366"
$values[0,10] = [double]"44.956116621581003"
$values[0,11] = [double]"8.736349453978158E-2"
$values[1,0] = [double]"63.398963730569939"
$values[1,1] = [double]"25.996092575894199"
$values[1,2] = [double]"17.757670306812273"
$values[1,3] = [double]"11.689119170984457"
$values[1,4] = [double]"40.381424706943193"
$values[1,5] = [double]"24.434737389495581"
$values[1,6] = [double]"27.924006908462871"
$values[1,7] = [double]"18.950706342049894"
$values[1,8] = [double]"21.574622984919397"
$values[1,9] = [double]"32.483592400690846"
$values[1,10] = [double]"52.241959723474615"
$values[1,11] = [double]"1.2522100884035361"
$values[2,0] = [double]"28.390328151986182"
$values[2,1] = [double]"23.062518785692816"
$values[2,2] = [double]"41.553822152886113"
$values[2,3] = [double]"4.1934369602763386"
$values[2,4] = [double]"31.520589119326726"
$values[2,5] = [double]"37.695267810712416"
$values[2,6] = [double]"13.727115716753024"
$values[2,7] = [double]"15.488127442140067"
$values[2,8] = [double]"47.128445137805507"
$values[2,9] = [double]"14.21070811744387"
$values[2,10] = [double]"50.035767959122332"
$values[2,11] = [double]"2.7706708268330726"
$values[3,0] = [double]"10.756476683937825"
$values[3,1] = [double]"12.852119026149685"
$values[3,2] = [double]"70.238169526781078"
$values[3,3] = [double]"2.8981001727115712"
$values[3,4] = [double]"27.030057108506167"
$values[3,5] = [double]"45.349973998959946"
$values[3,6] = [double]"7.6476683937823822"
$values[3,7] = [double]"7.9377817853922457"
$values[3,8] = [double]"66.535621424856984"
$values[3,9] = [double]"11.326424870466321"
$values[3,10] = [double]"68.722272317403068"
$values[3,11] = [double]"6.4940197607904322"
$values[4,0] = [double]"5.6960276338514682"
$values[4,1] = [double]"6.6453261196272919"
$values[4,2] = [double]"77.924076963078519"
$values[4,3] = [double]"0.67012089810017261"
$values[4,4] = [double]"28.328524195972349"
$values[4,5] = [double]"51.673426937077473"
$values[4,6] = [double]"4.3834196891191706"
$values[4,7] = [double]"5.8097385031559972"
$values[4,8] = [double]"70.456578263130538"
$values[4,9] = [double]"4.9706390328151979"
$values[4,10] = [double]"52.109708446047499"
$values[4,11] = [double]"9.0629225169006755"
$values[5,0] = [double]"4.0725388601036263"
$values[5,1] = [double]"5.893898406973249"
$values[5,2] = [double]"84.715548621944876"
$values[5,3] = [double]"-0.84974093264248696"
$values[5,4] = [double]"26.30267508265705"
$values[5,5] = [double]"48.854914196567854"
$values[5,6] = [double]"3.1053540587219337"
$values[5,7] = [double]"5.0042079951908622"
$values[5,8] = [double]"64.580343213728554"
$values[5,9] = [double]"7.3886010362694305"
$values[5,10] = [double]"51.370303576795912"
$values[5,11] = [double]"12.328653146125845"
$values[6,0] = [double]"1.4473229706390329"
$values[6,1] = [double]"5.4851217312894498"
$values[6,2] = [double]"85.98439937597503"
$values[6,3] = [double]"-1.0051813471502591"
$values[6,4] = [double]"25.629395852119028"
$values[6,5] = [double]"53.597503900155999"
$values[6,6] = [double]"1.2918825561312606"
$values[6,7] = [double]"5.9660354673880365"
$values[6,8] = [double]"75.542381695267807"
$values[6,9] = [double]"3.6407599309153711"
$values[6,10] = [double]"47.955816050495947"
$values[6,11] = [double]"13.051482059282371"
$ws.Range("B2:M8").Value2 = $values

# ---------------------------------------------------------------------------
# 4) Update the sheet selection to match where the author ended up.
# ---------------------------------------------------------------------------
$ws.Range("D2:D8").Select() | Out-Null

Write-Output "edit complete"
